# Update packet-counter snapshot values for sheets R1, R3, SW1 and append
# the newly-observed interfaces (Loopback0 on R1, Vlan10 on SW1).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet R1
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("R1")

# Row 3 - Ethernet1/0
$ws.Cells.Item(3, 6).Value = 57118    # F3 rx_octets
$ws.Cells.Item(3, 7).Value = 620      # G3 rx_unicast
$ws.Cells.Item(3, 10).Value = 74050   # J3 tx_octets
$ws.Cells.Item(3, 11).Value = 470     # K3 tx_unicast

# Row 11 - FastEthernet0/0
$ws.Cells.Item(11, 2).Value = 378     # B11 rx_broadcast
$ws.Cells.Item(11, 6).Value = 28222   # F11 rx_octets
$ws.Cells.Item(11, 7).Value = 378     # G11 rx_unicast
$ws.Cells.Item(11, 10).Value = 8760   # J11 tx_octets
$ws.Cells.Item(11, 11).Value = 71     # K11 tx_unicast

# Row 12 - FastEthernet0/1
$ws.Cells.Item(12, 10).Value = 8400   # J12 tx_octets
$ws.Cells.Item(12, 11).Value = 69     # K12 tx_unicast

# Row 13 - new interface Loopback0
$ws.Cells.Item(13, 1).Value = "Loopback0"
$ws.Cells.Item(13, 2).Value = 0
$ws.Cells.Item(13, 3).Value = 0
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(13, 5).Value = 0
$ws.Cells.Item(13, 6).Value = 0
$ws.Cells.Item(13, 7).Value = 0
$ws.Cells.Item(13, 8).Value = 0
$ws.Cells.Item(13, 9).Value = 0
$ws.Cells.Item(13, 10).Value = 456
$ws.Cells.Item(13, 11).Value = 6
$ws.Range("A13:K13").HorizontalAlignment = -4108
$ws.Range("A13:K13").VerticalAlignment = -4108

# ---------------------------------------------------------------------
# Sheet R3
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("R3")

# Row 3 - em0
$ws.Cells.Item(3, 6).Value = 262608   # F3 rx_octets
$ws.Cells.Item(3, 10).Value = 562965  # J3 tx_octets

# Row 4 - em1
$ws.Cells.Item(4, 10).Value = 682     # J4 tx_octets

# Row 5 - em2
$ws.Cells.Item(5, 10).Value = 640     # J5 tx_octets

# ---------------------------------------------------------------------
# Sheet SW1
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("SW1")

# Row 3 - GigabitEthernet0/0
$ws.Cells.Item(3, 2).Value = 30       # B3 rx_broadcast
$ws.Cells.Item(3, 6).Value = 30533    # F3 rx_octets
$ws.Cells.Item(3, 7).Value = 342      # G3 rx_unicast
$ws.Cells.Item(3, 10).Value = 94903   # J3 tx_octets
$ws.Cells.Item(3, 11).Value = 523     # K3 tx_unicast

# Row 4 - GigabitEthernet0/1
$ws.Cells.Item(4, 2).Value = 11       # B4 rx_broadcast
$ws.Cells.Item(4, 5).Value = 1        # E4 rx_multicast
$ws.Cells.Item(4, 6).Value = 6324     # F4 rx_octets
$ws.Cells.Item(4, 7).Value = 53       # G4 rx_unicast
$ws.Cells.Item(4, 10).Value = 32811   # J4 tx_octets
$ws.Cells.Item(4, 11).Value = 400     # K4 tx_unicast

# Row 5 - GigabitEthernet0/2
$ws.Cells.Item(5, 10).Value = 33047   # J5 tx_octets
$ws.Cells.Item(5, 11).Value = 402     # K5 tx_unicast

# Row 6 - GigabitEthernet0/3
$ws.Cells.Item(6, 10).Value = 33047   # J6 tx_octets
$ws.Cells.Item(6, 11).Value = 402     # K6 tx_unicast

# Row 9 - GigabitEthernet1/2
$ws.Cells.Item(9, 10).Value = 427     # J9 tx_octets
$ws.Cells.Item(9, 11).Value = 1       # K9 tx_unicast

# Row 18 - new interface Vlan10
$ws.Cells.Item(18, 1).Value = "Vlan10"
$ws.Cells.Item(18, 2).Value = 0
$ws.Cells.Item(18, 3).Value = 0
$ws.Cells.Item(18, 4).Value = 0
$ws.Cells.Item(18, 5).Value = 0
$ws.Cells.Item(18, 6).Value = 0
$ws.Cells.Item(18, 7).Value = 0
$ws.Cells.Item(18, 8).Value = 0
$ws.Cells.Item(18, 9).Value = 0
$ws.Cells.Item(18, 10).Value = 0
$ws.Cells.Item(18, 11).Value = 0
$ws.Range("A18:K18").HorizontalAlignment = -4108
$ws.Range("A18:K18").VerticalAlignment = -4108
